$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.692.61"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").Value = "2.732.03"
$ws.Range("E3").Value = "  -0.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.93"
$ws.Range("E5").Value = "  -1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.15"
$ws.Range("E6").Value = "  +2.18%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -0.47%  "

# Row 9
$ws.Range("E9").Value = "  +0.58%  "

# Row 10
$ws.Range("E10").Value = "  +3.82%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  +4.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.378"
$ws.Range("E12").Value = "  +0.00%  "

# Row 13
$ws.Range("D13").Value = "3.218.51"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.86"
$ws.Range("E14").Value = "  +2.32%  "

# Row 15
$ws.Range("D15").Value = "63.548.98"
$ws.Range("E15").Value = "  +0.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  +1.05%  "

# Row 17
$ws.Range("D17").Value = "2.739.37"
$ws.Range("E17").Value = "  -0.07%  "

# Row 18
$ws.Range("E18").Value = "  +3.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  -0.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.50"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("E21").Value = "  -1.88%  "

# Row 22
$ws.Range("E22").Value = "  +0.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.521"
$ws.Range("E23").Value = "  -2.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.10"
$ws.Range("E24").Value = "  -1.09%  "

# Row 25
$ws.Range("E25").Value = "  +0.77%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.35"
$ws.Range("E27").Value = "  +0.54%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0903"
$ws.Range("E28").Value = "  +3.06%  "

# Row 29
$ws.Range("E29").Value = "  +2.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  +4.96%  "

# Row 31
$ws.Range("E31").Value = "  +12.33%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.70"
$ws.Range("E32").Value = "  -2.81%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.04"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("E34").Value = "  +1.46%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  +3.79%  "

# Row 37
$ws.Range("E37").Value = "  +1.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.975"
$ws.Range("E38").Value = "  +1.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "344.65"
$ws.Range("E39").Value = "  +7.13%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.25"
$ws.Range("E40").Value = "  +1.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.08"
$ws.Range("E41").Value = "  +0.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.38"
$ws.Range("E42").Value = "  -0.63%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.81"
$ws.Range("E43").Value = "  +3.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.03"
$ws.Range("E44").Value = "  -0.07%  "

# Row 45
$ws.Range("E45").Value = "  +0.40%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.85"
$ws.Range("E46").Value = "  -0.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.621"
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("E48").Value = "  -0.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("E49").Value = "  -0.02%  "

# Row 50
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.04"
$ws.Range("E51").Value = "  -0.10%  "
